# Auto-generated edit script: updates cryptos list prices/volumes (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '25.959.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").Value = "'" + '1.634.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.08%  '
$ws.Range("D4").Value = "'" + '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'" + '215.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("E6").Value = '  -1.57%  '
$ws.Range("D7").Value = "'" + '1.009'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = "'" + '0.2565'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = "'" + '0.06386'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = "'" + '19.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.57%  '
$ws.Range("D11").Value = "'" + '0.07740'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").Value = "'" + '1.635.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("D13").Value = "'" + '4.229'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.66%  '
$ws.Range("D14").Value = "'" + '1.861.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").Value = "'" + '0.5428'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'" + '0.0' + [char]0x2085 + '7905'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").Value = "'" + '63.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").Value = "'" + '25.979.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("D19").Value = "'" + '1.009'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = "'" + '204.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.84%  '
$ws.Range("D21").Value = "'" + '4.330'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").Value = "'" + '9.955'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").Value = "'" + '5.945'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").Value = "'" + '1.009'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = "'" + '1.962'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.82%  '
$ws.Range("D26").Value = "'" + '141.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("D27").Value = "'" + '0.1153'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = "'" + '6.816'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = "'" + '15.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'" + '1.238'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = "'" + '0.04994'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.43%  '
$ws.Range("D32").Value = "'" + '3.256'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").Value = "'" + '3.181'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("D34").Value = "'" + '1.533'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.02%  '
$ws.Range("E35").Value = '  -2.19%  '
$ws.Range("D36").Value = "'" + '0.9009'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.44%  '
$ws.Range("D37").Value = "'" + '2.642'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.19%  '
$ws.Range("D38").Value = "'" + '0.5645'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("D39").Value = "'" + '1.125.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").Value = "'" + '0.01554'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("D41").Value = "'" + '1.009'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").Value = "'" + '2.557'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = "'" + '5.634'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").Value = "'" + '0.8106'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.15%  '
$ws.Range("D45").Value = "'" + '99.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("D46").Value = "'" + '1.772.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("D47").Value = "'" + '0.0' + [char]0x2088 + '114'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.61%  '
$ws.Range("D48").Value = "'" + '0.4527'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").Value = "'" + '54.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'" + '0.05050'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.85%  '
